$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1002.8
$ws.Range("I8").Value = 1002.8
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3008.4
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2869.4
$ws.Range("N8").ClearContents()
$ws.Range("H10").Value = 1756.75
$ws.Range("J10").Value = 3011.5
$ws.Range("L10").Value = 3011.5
$ws.Range("N10").Value = -3597.5
$ws.Range("H33").Value = 1015.7273
$ws.Range("I33").Value = 299.3125
$ws.Range("K33").Value = 299.3125
$ws.Range("M33").Value = -70.3125
$ws.Range("H100").Value = 1989.6154
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082
$ws.Range("H116").Value = 7037.3125
$ws.Range("I116").Value = 6207.0586
$ws.Range("K116").Value = 6207.0586
$ws.Range("M116").Value = -2765.0586
$ws.Range("H137").Value = 1182.8334
$ws.Range("I137").Value = 999.2857
$ws.Range("K137").Value = 2997.8571
$ws.Range("M137").Value = -447.8571000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 5722.273
$ws.Range("I32").Value = 5690.7334
$ws.Range("K32").Value = 5690.7334
$ws.Range("M32").Value = -5403.7334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3977.3076
$ws.Range("I134").Value = 3892.0833
$ws.Range("K134").Value = 11676.2499
$ws.Range("M134").Value = -9141.249899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2106.0725
$ws.Range("I31").Value = 1089.0227
$ws.Range("K31").Value = 1089.0227
$ws.Range("M31").Value = -794.0227
$ws.Range("H34").Value = 2106.0725
$ws.Range("I34").Value = 1089.0227
$ws.Range("K34").Value = 1089.0227
$ws.Range("M34").Value = -887.0227
$ws.Range("H62").Value = 127590.75
$ws.Range("I62").Value = 252134
$ws.Range("K62").Value = 252134
$ws.Range("M62").Value = -251510
$ws.Range("H65").Value = 127590.75
$ws.Range("I65").Value = 252134
$ws.Range("K65").Value = 1260670
$ws.Range("M65").Value = -1257550
$ws.Range("H86").Value = 11518.889
$ws.Range("I86").Value = 6238
$ws.Range("J86").Value = 13027.714
$ws.Range("K86").Value = 6238
$ws.Range("L86").Value = 13027.714
$ws.Range("M86").Value = -5115
$ws.Range("N86").Value = -15273.714
$ws.Range("H89").Value = 11518.889
$ws.Range("I89").Value = 6238
$ws.Range("J89").Value = 13027.714
$ws.Range("K89").Value = 31190
$ws.Range("L89").Value = 65138.57
$ws.Range("M89").Value = -25574
$ws.Range("N89").Value = -76370.57000000001
$ws.Range("H107").Value = 2049.8235
$ws.Range("I107").Value = 1496.2858
$ws.Range("J107").Value = 4633
$ws.Range("K107").Value = 1496.2858
$ws.Range("L107").Value = 4633
$ws.Range("M107").Value = 423.7141999999999
$ws.Range("N107").Value = -8473
$ws.Range("H132").Value = 2444.5454
$ws.Range("I132").Value = 2100.111
$ws.Range("K132").Value = 6300.333
$ws.Range("M132").Value = -3770.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 412.66666
$ws.Range("I7").Value = 369.16666
$ws.Range("J7").Value = 499.66666
$ws.Range("K7").Value = 1107.49998
$ws.Range("L7").Value = 1498.99998
$ws.Range("M7").Value = -995.4999800000001
$ws.Range("N7").Value = -1722.99998
$ws.Range("H25").Value = 4000.7778
$ws.Range("I25").Value = 1999.6666
$ws.Range("K25").Value = 5998.9998
$ws.Range("M25").Value = -5829.9998
$ws.Range("H30").Value = 4000.7778
$ws.Range("I30").Value = 1999.6666
$ws.Range("K30").Value = 5998.9998
$ws.Range("M30").Value = -5896.9998
$ws.Range("H41").Value = 200
$ws.Range("H44").Value = 2006.1538
$ws.Range("I44").Value = 177.09091
$ws.Range("J44").Value = 3347.4666
$ws.Range("K44").Value = 531.27273
$ws.Range("L44").Value = 10042.3998
$ws.Range("M44").Value = -133.27273
$ws.Range("N44").Value = -10838.3998
$ws.Range("H120").Value = 14022
$ws.Range("J120").Value = 16033
$ws.Range("L120").Value = 48099
$ws.Range("N120").Value = -57775

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 660.1905
$ws.Range("I2").Value = 206.58824
$ws.Range("K2").Value = 206.58824
$ws.Range("M2").Value = -93.58824000000001
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H97").Value = 50001556
$ws.Range("I97").Value = 66668284
$ws.Range("K97").Value = 66668284
$ws.Range("M97").Value = -66667788
$ws.Range("H107").Value = 440.9091
$ws.Range("I107").Value = 375.3125
$ws.Range("J107").Value = 615.8333
$ws.Range("K107").Value = 375.3125
$ws.Range("L107").Value = 615.8333
$ws.Range("M107").Value = 1544.6875
$ws.Range("N107").Value = -4455.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 33855
$ws.Range("I43").Value = 32011
$ws.Range("J43").Value = 34162.332
$ws.Range("K43").Value = 32011
$ws.Range("L43").Value = 34162.332
$ws.Range("M43").Value = -31818
$ws.Range("N43").Value = -34548.332
$ws.Range("H46").Value = 2615.2778
$ws.Range("I46").Value = 673.1667
$ws.Range("K46").Value = 673.1667
$ws.Range("M46").Value = -485.1667
$ws.Range("H75").Value = 175000
$ws.Range("I75").Value = 175000
$ws.Range("J75").Value = 175000
$ws.Range("K75").Value = 175000
$ws.Range("L75").Value = 175000
$ws.Range("M75").Value = -174064
$ws.Range("N75").Value = -176872
$ws.Range("H78").Value = 175000
$ws.Range("I78").Value = 175000
$ws.Range("J78").Value = 175000
$ws.Range("K78").Value = 525000
$ws.Range("L78").Value = 525000
$ws.Range("M78").Value = -520320
$ws.Range("N78").Value = -534360
$ws.Range("H82").Value = 1549.3793
$ws.Range("I82").Value = 1555.7391
$ws.Range("J82").Value = 1525
$ws.Range("K82").Value = 1555.7391
$ws.Range("L82").Value = 1525
$ws.Range("M82").Value = -1194.7391
$ws.Range("N82").Value = -2247
$ws.Range("H85").Value = 1549.3793
$ws.Range("I85").Value = 1555.7391
$ws.Range("J85").Value = 1525
$ws.Range("K85").Value = 1555.7391
$ws.Range("L85").Value = 1525
$ws.Range("M85").Value = -307.7391
$ws.Range("N85").Value = -4021
$ws.Range("H101").Value = 38250.375
$ws.Range("J101").Value = 38250.375
$ws.Range("L101").Value = 38250.375
$ws.Range("N101").Value = -44740.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18265.834
$ws.Range("I62").Value = 18108.182
$ws.Range("K62").Value = 18108.182
$ws.Range("M62").Value = -17484.182
$ws.Range("H65").Value = 18265.834
$ws.Range("I65").Value = 18108.182
$ws.Range("K65").Value = 90540.91
$ws.Range("M65").Value = -87420.91
$ws.Range("H113").Value = 870.9355
$ws.Range("J113").Value = 545.375
$ws.Range("L113").Value = 1636.125
$ws.Range("N113").Value = -5976.125
$ws.Range("H126").Value = 33027504
$ws.Range("I126").Value = 40953830
$ws.Range("J126").Value = 1137.8334
$ws.Range("K126").Value = 122861490
$ws.Range("L126").Value = 3413.5002
$ws.Range("M126").Value = -122859020
